# Apply "Allow for different key viewer profiles" translation updates.
$wb = $excel.ActiveWorkbook

# --- KeyLimiter sheet: add trailing period to Spanish description ---
$wsLimiter = $wb.Worksheets.Item("KeyLimiter")
$wsLimiter.Range("D3").Value = "Cuando está activo, restringe qué teclas cuentan como válidas."

# --- KeyViewer sheet: add trailing period to DESCRIPTION row translations ---
$wsViewer = $wb.Worksheets.Item("KeyViewer")
$wsViewer.Range("B3").Value = "Shows a key viewer for registered keys."
$wsViewer.Range("C3").Value = "등록된 키들의 키뷰어 보이기."
$wsViewer.Range("D3").Value = "Muestra teclas en pantalla."

# --- KeyViewer sheet: append new rows for the Profiles feature ---
$wsViewer.Range("A19").Value = "PROFILES"
$wsViewer.Range("B19").Value = "Profiles:"
$wsViewer.Range("C19").Value = "프로파일들:"
$wsViewer.Range("D19").Value = "Perfiles:"

$wsViewer.Range("A20").Value = "PROFILE_NAME"
$wsViewer.Range("B20").Value = "Profile name:"
$wsViewer.Range("C20").Value = "프로파일 이름:"
$wsViewer.Range("D20").Value = "Nombre del perfil:"

$wsViewer.Range("A21").Value = "NEW"
$wsViewer.Range("B21").Value = "New"
$wsViewer.Range("C21").Value = "생성"
$wsViewer.Range("D21").Value = "Nuevo"

$wsViewer.Range("A22").Value = "DUPLICATE"
$wsViewer.Range("B22").Value = "Duplicate"
$wsViewer.Range("C22").Value = "복제"
$wsViewer.Range("D22").Value = "Duplicar"

$wsViewer.Range("A23").Value = "DELETE"
$wsViewer.Range("B23").Value = "Delete"
$wsViewer.Range("C23").Value = "삭제"
$wsViewer.Range("D23").Value = "Eliminar"
